$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = @{ D = "41.922.89";  E = "  -2.11%  " }
    3  = @{ D = "2.488.47";   E = "  -3.56%  " }
    4  = @{ D = "0.999";      E = "  -0.06%  " }
    5  = @{ D = "301.14";     E = "  -0.45%  " }
    6  = @{ D = "94.18";      E = "  -2.82%  " }
    7  = @{ E = "  +0.38%  " }
    8  = @{ E = "  +0.12%  " }
    9  = @{ D = "0.526";      E = "  -4.15%  " }
    10 = @{ D = "35.84";      E = "  -2.55%  " }
    11 = @{ E = "  -1.01%  " }
    12 = @{ D = "0.111";      E = "  -2.58%  " }
    13 = @{ D = "7.34";       E = "  -4.92%  " }
    14 = @{ D = "2.866.74";   E = "  -3.81%  " }
    15 = @{ D = "2.491.34";   E = "  -4.38%  " }
    16 = @{ D = "14.81";      E = "  +3.17%  " }
    17 = @{ D = "0.846";      E = "  -4.35%  " }
    18 = @{ D = "41.964.07";  E = "  -2.13%  " }
    19 = @{ D = "12.70";      E = "  -1.41%  " }
    20 = @{ D = "0.0₃0959";   E = "  -3.03%  " }
    21 = @{ D = "6.33";       E = "  -4.66%  " }
    22 = @{ D = "70.55";      E = "  -1.86%  " }
    23 = @{ D = "246.75";     E = "  -3.03%  " }
    24 = @{ E = "  -2.56%  " }
    25 = @{ D = "1.98";       E = "  -6.60%  " }
    26 = @{ D = "26.44";      E = "  -7.63%  " }
    27 = @{ D = "0.997";      E = "  -0.35%  " }
    28 = @{ D = "2.29";       E = "  +8.29%  " }
    29 = @{ D = "10.02";      E = "  -1.75%  " }
    30 = @{ D = "36.86";      E = "  -6.25%  " }
    31 = @{ D = "5.82";       E = "  -3.36%  " }
    32 = @{ D = "153.14";     E = "  -1.33%  " }
    33 = @{ E = "  -3.14%  " }
    34 = @{ E = "  -5.45%  " }
    35 = @{ D = "0.0772";     E = "  -5.17%  " }
    36 = @{ E = "  -6.26%  " }
    37 = @{ D = "18.17";      E = "  -1.05%  " }
    38 = @{ E = "  -1.38%  " }
    39 = @{ E = "  -1.51%  " }
    40 = @{ D = "23.82";      E = "  +1.83%  " }
    41 = @{ D = "3.79";       E = "  -2.49%  " }
    42 = @{ D = "3.32";       E = "  -2.40%  " }
    43 = @{ D = "0.998";      E = "  -0.13%  " }
    44 = @{ D = "2.035.52";   E = "  -1.78%  " }
    45 = @{ D = "0.0294";     E = "  -4.88%  " }
    46 = @{ D = "1.91";       E = "  -8.53%  " }
    47 = @{ D = "8.85";       E = "  -4.25%  " }
    48 = @{ D = "82.61";      E = "  -3.00%  " }
    49 = @{ D = "2.727.52";   E = "  -3.80%  " }
    50 = @{ D = "71.13";      E = "  -6.21%  " }
    51 = @{ D = "0.185";      E = "  -2.65%  " }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    if ($cells.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cells["D"]
        $cell.ClearFormats()
    }
    if ($cells.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cells["E"]
        $cell.ClearFormats()
    }
}
